$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("...nights", "～泊|～はく")
  ,@("with...", "～付|～つき")
  ,@("with breakfast", "朝食付き|ちょうしょくつき")
  ,@("one night with two meals", "一泊二食付|いっぱくにしょくつき")
  ,@("checking in", "チェックイン（する）")
  ,@("checking out", "チェックアウト（する）")
  ,@("single room", "シングル")
  ,@("double room", "ダブル")
  ,@("twin room", "ツイン")
  ,@("...person(s)", "～名|～めい")
  ,@("receptionist; front desk", "フロント")
  ,@("non-smoking", "禁煙|きんえん")
  ,@("smoking", "喫煙|きつえん")
  ,@("large spa", "大浴場|だいよくじょう")
  ,@("I have a reservation under the name of XX.", "予約した○○です。|よやくしたまるまるです。")
  ,@("Could you check me in, please?", "チェックインお願いします。|チェックインおねがいします。")
  ,@("Can I pay by XX?", "○○で払えますか。|まるまるではらえますか。")
  ,@("Can you change US dollars to Japanese yen?", "アメリカドルを日本円に両替できますか。|アメリカドルをにほんえんにりょうがえできますか。")
  ,@("Would you call a taxi for us?", "タクシーを呼んでもらえますか。|タクシーをよんでもらえますか。")
  ,@("Are there any restaurants you recommend near here?", "近くにおすすめのレストランがありますか。|ちかくにおすすめのレストランがありますか。")
  ,@("Could you keep my luggage until 2 o'clock?", "二時まで荷物を預かってくれませんか。|にじまでにもつをあずかってくれませんか。")
)

$startRow = 68
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
